# Update "想去人数" (F) / "最低票价" (G) counts that changed between scrapes,
# across the four worksheets of the workbook.
$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 421
$ws.Range("F10").Value = 394
$ws.Range("F11").Value = 439
$ws.Range("F12").Value = 35
$ws.Range("F14").Value = 373
$ws.Range("F15").Value = 56
$ws.Range("F16").Value = 68
$ws.Range("F17").Value = 21
$ws.Range("F18").Value = 600
$ws.Range("F19").Value = 1489
$ws.Range("F20").Value = 5799
$ws.Range("G20").Value = 68
$ws.Range("F22").Value = 1631
$ws.Range("F26").Value = 5464
$ws.Range("F27").Value = 135
$ws.Range("F29").Value = 1564
$ws.Range("F30").Value = 621
$ws.Range("F32").Value = 74
$ws.Range("F33").Value = 1103
$ws.Range("F36").Value = 10
$ws.Range("F38").Value = 3820

# Sheet: 演出 (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 80
$ws.Range("F5").Value = 185
$ws.Range("F8").Value = 272

# Sheet: 本地生活 (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 9438
$ws.Range("F3").Value = 588
$ws.Range("F4").Value = 2171
$ws.Range("F5").Value = 417

# Sheet: 全部类型 (All types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 9438
$ws.Range("F3").Value = 588
$ws.Range("F4").Value = 2171
$ws.Range("F6").Value = 421
$ws.Range("F11").Value = 394
$ws.Range("F12").Value = 439
$ws.Range("F14").Value = 373
$ws.Range("F15").Value = 56
$ws.Range("F16").Value = 68
$ws.Range("F19").Value = 1489
$ws.Range("F20").Value = 5799
$ws.Range("G20").Value = 68
$ws.Range("F22").Value = 1631
$ws.Range("F28").Value = 5464
$ws.Range("F29").Value = 135
$ws.Range("F31").Value = 1564
$ws.Range("F32").Value = 622
$ws.Range("F34").Value = 1103
$ws.Range("F40").Value = 10
$ws.Range("F45").Value = 3820
